# Apply the commit's edits to the active document.
$d = $word.ActiveDocument

# --- 1. document.xml: give the (only, currently empty) paragraph the
#        "Block Text" style and insert the run of text "asa" before the
#        existing _GoBack bookmark markers. ------------------------------
$p = $d.Paragraphs(1)
$r = $d.Range(0, 0)
$r.InsertBefore("asa")
$p.Style = "Block Text"

# --- 2. styles.xml: tweak the BlockText paragraph style definition -------
$s = $d.Styles("Block Text")

# <w:ind w:left="720"/> inside <w:pPr> (720 twips = 36 pt = 0.5")
$s.ParagraphFormat.LeftIndent = 36

# <w:i/> inside <w:rPr>
$s.Font.Italic = $true

# <w:lang w:val="en-CA"/> inside <w:rPr>
$s.LanguageID = "en-CA"
